$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "43.005.07"
Set-TextValue $ws "E2" "  +4.29%  "
Set-TextValue $ws "D3" "2.240.12"
Set-TextValue $ws "E3" "  +3.12%  "
Set-TextValue $ws "E4" "  +0.10%  "
Set-TextValue $ws "D5" "245.95"
Set-TextValue $ws "E5" "  +3.92%  "
Set-TextValue $ws "D6" "0.620"
Set-TextValue $ws "E6" "  +0.68%  "
Set-TextValue $ws "D7" "76.34"
Set-TextValue $ws "E7" "  +8.45%  "
Set-TextValue $ws "E8" "  -0.05%  "
Set-TextValue $ws "D9" "0.616"
Set-TextValue $ws "E9" "  +6.11%  "
Set-TextValue $ws "D10" "40.94"
Set-TextValue $ws "E10" "  +1.45%  "
Set-TextValue $ws "D11" "0.0934"
Set-TextValue $ws "E11" "  +0.43%  "
Set-TextValue $ws "D12" "55.57"
Set-TextValue $ws "E12" "  +1.09%  "
Set-TextValue $ws "E13" "  +2.38%  "
Set-TextValue $ws "E14" "  +0.59%  "
Set-TextValue $ws "D15" "2.553.97"
Set-TextValue $ws "E15" "  +2.35%  "
Set-TextValue $ws "D16" "14.65"
Set-TextValue $ws "E16" "  +5.16%  "
Set-TextValue $ws "D17" "2.253.08"
Set-TextValue $ws "E17" "  +3.82%  "
Set-TextValue $ws "D18" "0.812"
Set-TextValue $ws "E18" "  +0.55%  "
Set-TextValue $ws "D19" "42.929.27"
Set-TextValue $ws "E19" "  +4.73%  "
Set-TextValue $ws "D20" "0.0000105"
Set-TextValue $ws "E20" "  +2.81%  "
Set-TextValue $ws "D21" "71.24"
Set-TextValue $ws "E21" "  +0.99%  "
Set-TextValue $ws "E22" "  +0.79%  "
Set-TextValue $ws "D23" "10.25"
Set-TextValue $ws "E23" "  +3.91%  "
Set-TextValue $ws "D24" "2.23"
Set-TextValue $ws "E24" "  +13.74%  "
Set-TextValue $ws "D25" "230.61"
Set-TextValue $ws "E25" "  +1.74%  "
Set-TextValue $ws "E26" "  +0.01%  "
Set-TextValue $ws "E27" "  -0.03%  "
Set-TextValue $ws "E28" "  -5.64%  "
Set-TextValue $ws "E29" "  +2.28%  "
Set-TextValue $ws "D30" "38.27"
Set-TextValue $ws "E30" "  +25.32%  "
Set-TextValue $ws "D31" "174.20"
Set-TextValue $ws "E31" "  +3.77%  "
Set-TextValue $ws "E32" "  -2.28%  "
Set-TextValue $ws "D33" "20.30"
Set-TextValue $ws "E33" "  +1.46%  "
Set-TextValue $ws "D34" "0.0797"
Set-TextValue $ws "E34" "  +3.38%  "
Set-TextValue $ws "D35" "5.37"
Set-TextValue $ws "E35" "  +3.75%  "
Set-TextValue $ws "D36" "0.114"
Set-TextValue $ws "E36" "  +10.96%  "
Set-TextValue $ws "D37" "0.123"
Set-TextValue $ws "E37" "  +1.06%  "
Set-TextValue $ws "D38" "4.39"
Set-TextValue $ws "E38" "  +6.59%  "
Set-TextValue $ws "E39" "  +14.31%  "
Set-TextValue $ws "D40" "13.02"
Set-TextValue $ws "E40" "  +8.31%  "
Set-TextValue $ws "E41" "  +2.83%  "
Set-TextValue $ws "E42" "  +2.31%  "
Set-TextValue $ws "D43" "0.204"
Set-TextValue $ws "E43" "  +6.92%  "
Set-TextValue $ws "D44" "60.06"
Set-TextValue $ws "E44" "  -0.27%  "
Set-TextValue $ws "D45" "105.74"
Set-TextValue $ws "E45" "  +7.72%  "
Set-TextValue $ws "D46" "8.66"
Set-TextValue $ws "E46" "  +4.05%  "
Set-TextValue $ws "D47" "0.0990"
Set-TextValue $ws "E47" "  +1.38%  "
Set-TextValue $ws "B48" "ARBITRUM"
Set-TextValue $ws "C48" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws "D48" "1.10"
Set-TextValue $ws "E48" "  +1.18%  "
Set-TextValue $ws "B49" "NEARProtocol"
Set-TextValue $ws "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws "D49" "2.31"
Set-TextValue $ws "E49" "  +4.00%  "
Set-TextValue $ws "B50" "WOONetwork"
Set-TextValue $ws "C50" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue $ws "D50" "0.439"
Set-TextValue $ws "E50" "  +15.49%  "
Set-TextValue $ws "D51" "1.15"
Set-TextValue $ws "E51" "  +0.91%  "
